$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.848.48"
$ws.Range("E2").Value = "  -0.20%  "

$ws.Range("D3").Value = "1.886.18"
$ws.Range("E3").Value = "  -0.52%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("E5").Value = "  -4.52%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "242.53"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.65%  "

$ws.Range("E7").Value = "  -0.15%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3116"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -0.89%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "25.45"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -1.09%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.07123"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -2.05%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.08472"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +4.37%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.7594"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -1.91%  "

$ws.Range("D13").Value = "1.907.37"
$ws.Range("E13").Value = "  +0.21%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "5.357"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -2.20%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "93.30"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -1.13%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "6.141"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -1.23%  "

$ws.Range("D17").Value = "29.902.66"
$ws.Range("E17").Value = "  +0.02%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "13.70"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -1.79%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "243.07"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -1.24%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "0.000007790"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.56%  "

$ws.Range("D21").Value = "2.156.15"
$ws.Range("E21").Value = "  +0.73%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.9993"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -0.29%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "8.010"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -1.63%  "

$ws.Range("E24").Value = "  -0.07%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "0.1591"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.82%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "9.377"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -0.95%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "162.65"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -1.10%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "18.74"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -0.15%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.028"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +0.27%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.510"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +5.34%  "

$ws.Range("E31").Value = "  -1.05%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "4.473"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -0.13%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "4.100"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +0.39%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.05396"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -3.33%  "

$ws.Range("E35").Value = "  -0.64%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.7434"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -1.45%  "

$ws.Range("E37").Value = "  +0.69%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "2.711"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +1.19%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.01933"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -0.16%  "

$ws.Range("E40").Value = "  -0.70%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.4454"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -0.15%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "6.076"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +1.85%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "72.65"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -1.91%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.089.76"
$ws.Range("E44").Value = "  -4.64%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.8600"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +0.68%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.9999"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -0.29%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "102.56"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +0.42%  "

$ws.Range("E48").Value = "  +1.83%  "

$ws.Range("E49").Value = "  -1.58%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "3.053"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -2.97%  "

$ws.Range("D51").Value = "2.056.66"
$ws.Range("E51").Value = "  +2.30%  "
